# Daily update commit: "Updated: st 07. 10. 2021"
# Revises several historical AgTests (F) / AgPosit (G) values and appends
# a new day of data (row 581) to the DailyStats sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to previously reported AgTests / AgPosit figures ---
$ws.Range("F548").Value = 17258

$ws.Range("F549").Value = 10875

$ws.Range("F551").Value = 17789

$ws.Range("F553").Value = 15539

$ws.Range("F554").Value = 17954
$ws.Range("G554").Value = 185

$ws.Range("F555").Value = 21616

$ws.Range("F558").Value = 24677

$ws.Range("F560").Value = 6103

$ws.Range("F561").Value = 24249

$ws.Range("F562").Value = 27145

$ws.Range("F565").Value = 28835

$ws.Range("F566").Value = 25994

$ws.Range("F567").Value = 23466

$ws.Range("F568").Value = 23938

$ws.Range("F569").Value = 32380
$ws.Range("G569").Value = 358

$ws.Range("F570").Value = 15215

$ws.Range("F571").Value = 14998
$ws.Range("G571").Value = 278

$ws.Range("F572").Value = 33349
$ws.Range("G572").Value = 596

$ws.Range("F573").Value = 26981
$ws.Range("G573").Value = 400

$ws.Range("F574").Value = 23337
$ws.Range("G574").Value = 348

$ws.Range("F575").Value = 25907
$ws.Range("G575").Value = 386

$ws.Range("F576").Value = 28348
$ws.Range("G576").Value = 419

$ws.Range("F577").Value = 14295
$ws.Range("G577").Value = 269

$ws.Range("F578").Value = 14829
$ws.Range("G578").Value = 323

$ws.Range("F579").Value = 31479
$ws.Range("G579").Value = 609

# --- Fill in AgTests / AgPosit for the last existing row (580), which
#     previously only had columns A-E populated ---
$ws.Range("F580").Value = 27490
$ws.Range("G580").Value = 489

# --- Append the new day's row (581) ---
$ws.Range("A581").Value = 44475
$ws.Range("B581").Value = 420924
$ws.Range("C581").Value = 11115
$ws.Range("D581").Value = 1451
$ws.Range("E581").Value = 12705
$ws.Range("F581").Value = 19169
$ws.Range("G581").Value = 340
